# Handback status report generation: refresh the timestamp columns that are
# recomputed each time the handback report is regenerated.
#
# Overview sheet: "Latest HO Xliff Generate Date" (col G) for the
#   44d2b0a4-... row.
# zh-cn sheet: "Correspond Handoff Datetime" (col H) and
#   "Correspond Handback DateTime" (col K) for the 44d2b0a4-... row.
# de-de sheet: "Correspond Handoff Datetime" (col H, shared with the
#   Overview sheet's value) and "Correspond Handback DateTime" (col K)
#   for the 44d2b0a4-... row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G3 and de-de!H3 previously shared the same text
# ("2016-09-06 15:33:45"); both move to the same new timestamp.
$wsOverview.Range("G3").Value = "2016-09-06 15:35:39"
$wsDeDe.Range("H3").Value = "2016-09-06 15:35:39"

# zh-cn!H3 (Correspond Handoff Datetime) and zh-cn!K3 (Correspond
# Handback DateTime).
$wsZhCn.Range("H3").Value = "2016-09-06 15:35:33"
$wsZhCn.Range("K3").Value = "2016-09-06 15:35:54"

# de-de!K3 (Correspond Handback DateTime).
$wsDeDe.Range("K3").Value = "2016-09-06 15:36:20"
